# Actualización tras jornada 1
# Insert a new "1X2" column before the existing "jornada" block (old column E),
# which shifts the old E:I block (jornada, RTDO L, LOCAL, VISITANTE, RTDO V)
# one column to the right (becoming F:J), then populate the new leading
# column (E) and append a new trailing "1X2.1" column (K) with the
# jornada-1 results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new column at E - shifts old E:I (jornada..RTDO V) to F:J.
#    The engine mirrors Excel's default "inherit format from the left"
#    behaviour, so the new E column picks up the bold/bordered header style.
$ws.Range("E1").EntireColumn.Insert()

# 2) New header cells.
$ws.Range("E1").Value2 = "1X2"
$ws.Range("G1").Value2 = "RTDO L.1"
$ws.Range("H1").Value2 = "LOCAL.1"
$ws.Range("I1").Value2 = "VISITANTE.1"
$ws.Range("J1").Value2 = "RTDO V.1"

# 3) Append the trailing K column; copy formats from the neighbouring
#    header cell (J1) so it matches the rest of the header row's style,
#    then set its text.
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)
$ws.Range("K1").Value2 = "1X2.1"
$excel.CutCopyMode = $false

# 4) Row 2 (Palop/Lope vs Diego->Puche/Coquina) results.
$ws.Range("A2").Value2 = 43
$ws.Range("D2").Value2 = 19
$ws.Range("E2").Value2 = 1
$ws.Range("G2").Value2 = 36
$ws.Range("H2").Value2 = "Puche"
$ws.Range("J2").Value2 = 34
$ws.Range("K2").Value2 = 1

# 5) Row 3 (Kero/Fale vs Kike/Armada) results.
$ws.Range("A3").Value2 = 50
$ws.Range("D3").Value2 = 36
$ws.Range("E3").Value2 = 1
$ws.Range("G3").Value2 = 31
$ws.Range("J3").Value2 = 49
$ws.Range("K3").Value2 = 2

# 6) Row 4 (Tony/Ruso vs Papu/Gonzo) results.
$ws.Range("A4").Value2 = 41
$ws.Range("D4").Value2 = 23
$ws.Range("E4").Value2 = 1
$ws.Range("G4").Value2 = 38
$ws.Range("J4").Value2 = 32
$ws.Range("K4").Value2 = 1

Write-Output ("Used range after edit: " + $ws.UsedRange.Address())
